# Issue 404: Convert 'Coordinator' to 'Instructor' (Stage 3 diagram tweak).
# The "CoordData" boxes in the Db-layer data-transfer-classes diagram are
# renamed to "InstructorData" to match the rest of the app's terminology.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "CoordData") {
                $shp.TextFrame.TextRange.Text = "InstructorData"
            }
        }
    }
}
